# Workbook/worksheets
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# ---------------------------------------------------------------
# Hoja1 (sheet1) - "Tomb" zone column B re-shuffled / a new value
# ---------------------------------------------------------------

# B3: Slime -> Enemy Walker (style/format unchanged)
$ws1.Range("B3").Value = "Enemy Walker"

# B4: Enemy Walker -> Enemy Walker Walls (style/format unchanged)
$ws1.Range("B4").Value = "Enemy Walker Walls"

# B5: Enemy Walker Walls -> Enemy Eater (style/format unchanged)
$ws1.Range("B5").Value = "Enemy Eater"

# B6: Enemy Eater -> Enemy Shoot, and its formatting now matches column C
# (fill/border change s="14" -> s="3"), so copy the format from C6 first.
$ws1.Range("C6").Copy()
$ws1.Range("B6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws1.Range("B6").Value = "Enemy Shoot"

# B13: was blank, now contains "Nube voladora" (style/format unchanged)
$ws1.Range("B13").Value = "Nube voladora"

# ---------------------------------------------------------------
# Hoja2 (sheet2) - new lookup list entry
# ---------------------------------------------------------------

# New row 25: a new catalog entry "Cadenas Colgantes"
$ws2.Range("A25").Value = "Cadenas Colgantes"

# ---------------------------------------------------------------
# View state: active cell / selection / scroll position
# ---------------------------------------------------------------

# Hoja2 window scroll + selection
[void]$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws2.Range("D34").Select()

# Hoja1 stays the visible/active tab, with a new selected cell
[void]$ws1.Activate()
[void]$ws1.Range("F7").Select()
